$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap values in columns A, Q, R and AC between row 5 and row 6.
$cols = @("A", "Q", "R", "AC")

foreach ($col in $cols) {
    $cell5 = $ws.Range("$col`5")
    $cell6 = $ws.Range("$col`6")
    $tmp = $cell5.Value2
    $cell5.Value = $cell6.Value2
    $cell6.Value = $tmp
}
